$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as text so numeric-looking strings
# (with trailing zeros / thousands-style dots) are preserved exactly,
# matching the original inline-string cell content.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.192.06'
$ws.Range("E2").Value = '  -2.15%  '
$ws.Range("D3").Value = '1.671.47'
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '217.63'
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").Value = '0.5114'
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '0.2656'
$ws.Range("E8").Value = '  +3.07%  '
$ws.Range("D9").Value = '0.06378'
$ws.Range("E9").Value = '  +3.11%  '
$ws.Range("D10").Value = '21.51'
$ws.Range("E10").Value = '  -2.22%  '
$ws.Range("D11").Value = '0.07394'
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").Value = '1.668.31'
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("D13").Value = '4.544'
$ws.Range("E13").Value = '  +1.48%  '
$ws.Range("D14").Value = '0.5833'
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.000008630'
$ws.Range("E15").Value = '  +5.73%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '64.54'
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '26.314.46'
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '4.957'
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '1.007'
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '10.88'
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '189.21'
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").Value = '6.212'
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '1.008'
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").Value = '144.27'
$ws.Range("E24").Value = '  +0.90%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '7.653'
$ws.Range("E25").Value = '  +1.42%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.1178'
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '15.64'
$ws.Range("E27").Value = '  +2.86%  '
$ws.Range("B28").Value = 'Hedera'
$ws.Range("C28").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D28").Value = '0.05988'
$ws.Range("E28").Value = '  +1.85%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.287'
$ws.Range("E29").Value = '  -3.76%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '1.326'
$ws.Range("E30").Value = '  -1.50%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '3.528'
$ws.Range("E31").Value = '  +2.05%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '3.525'
$ws.Range("E32").Value = '  +2.76%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").Value = '1.648'
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '1.016'
$ws.Range("E34").Value = '  +2.98%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.6032'
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.377'
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = '2.663'
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01617'
$ws.Range("E38").Value = '  +1.30%  '
$ws.Range("D39").Value = '6.092'
$ws.Range("E39").Value = '  +3.37%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.081.55'
$ws.Range("E40").Value = '  -1.17%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '0.8714'
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '1.011'
$ws.Range("E42").Value = '  +0.58%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '100.16'
$ws.Range("E43").Value = '  +3.02%  '
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '1.821.97'
$ws.Range("E44").Value = '  -1.43%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.00000000115'
$ws.Range("E45").Value = '  +10.89%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '56.31'
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '1.012'
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '8.079'
$ws.Range("E48").Value = '  +1.39%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.05220'
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '0.4299'
$ws.Range("E50").Value = '  -0.85%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '5.893'
$ws.Range("E51").Value = '  +2.70%  '
